$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.091.23'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '2.575.26'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("D9").Value = '2.575.21'
$ws.Range("E9").Value = '  +2.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.35%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").Value = '3.042.40'
$ws.Range("E14").Value = '  +6.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000182'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.21%  '
$ws.Range("D17").Value = '70.021.19'
$ws.Range("E17").Value = '  +2.90%  '
$ws.Range("D18").Value = '2.581.32'
$ws.Range("E18").Value = '  +2.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.57%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '367.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.91%  '
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").Value = '2.704.13'
$ws.Range("E28").Value = '  +2.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").Value = '0.0₃0929'
$ws.Range("E30").Value = '  +0.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '520.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("E34").Value = '  +1.59%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.119'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.37%  '
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.97'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.15%  '
$ws.Range("E44").Value = '  -1.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '153.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.526'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("D50").Value = '0.0₆0261'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("E51").Value = '  +1.73%  '
